# Apply "updates done upto 17Jan2023" changes to TD.xlsx
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "HRM": consolidate the Admin/admin123 row up into row 4
# (removing the er6453/tryt76575 row entirely), then record a
# Pass/Fail verdict for every login attempt in column C.
# ---------------------------------------------------------------
$hrm = $wb.Worksheets.Item("HRM")

$hrm.Range("A4").Value = "Admin"
$hrm.Range("B4").Value = "admin123"
$hrm.Range("A5").Value = "gu5674"
$hrm.Range("B5").Value = "cvbncvnvcnv"
$hrm.Range("A6").Value = "gu5674sdgdsg"
$hrm.Range("B6").Value = "cvbncvnvcnvsdg"
$hrm.Rows.Item(7).Delete()

$hrm.Range("C2").Value = "Fail"
$hrm.Range("C3").Value = "Fail"
$hrm.Range("C4").Value = "Pass"
$hrm.Range("C5").Value = "Fail"
$hrm.Range("C6").Value = "Fail"
$hrm.Range("C2:C6").ClearFormats()

$hrm.Range("C2:C6").Select()

# ---------------------------------------------------------------
# Sheet "FBLoin": record a Pass/Fail verdict for every login
# attempt in column C (last row is the only valid credential).
# ---------------------------------------------------------------
$fbloin = $wb.Worksheets.Item("FBLoin")

$fbloin.Range("C2").Value = "Fail"
$fbloin.Range("C3").Value = "Fail"
$fbloin.Range("C4").Value = "Fail"
$fbloin.Range("C5").Value = "Fail"
$fbloin.Range("C6").Value = "Fail"
$fbloin.Range("C7").Value = "Fail"
$fbloin.Range("C8").Value = "Pass"

$fbloin.Activate()
$fbloin.Range("C2:C9").Select()
